# Mise à jour des résultats du script
# Adds a new data row (row 49) to the worksheet, mirroring the existing
# rows' layout: Date | Terme | Numéro de page | Occurences

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 49

# Keep column A as plain text (matches existing rows which store the date
# as a text string, not a native date value).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-04-14"

$ws.Cells.Item($row, 2).Value = "Rien ne nous concerne aujourd'hui !"

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "NA"

$ws.Cells.Item($row, 4).Value = 1
